# ---------------------------------------------------------------------------
# Insert a new "2022-Q4" sheet (fund-holding detail) right after "总计",
# and add a corresponding summary row on "总计" itself.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4, shifting the
#    existing three rows down by one, and renumber the index column (A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift existing data rows 2..4 down to 3..5 (bottom-up so nothing is
# clobbered before it is copied). Copy also carries the existing cell
# style (s="2" on column A) along with it.
$summary.Range("A4:D4").Copy($summary.Range("A5:D5"))
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

# New first data row: 2022-Q4
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 13
$summary.Range("D2").Value = 0.72

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q4" worksheet right after "总计". Duplicating the
#    "2021-Q4" sheet (10 rows) gives us the right sheetPr/pageMargins/
#    styles for free; we then overwrite its data and grow it to 14 rows.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item(4)
$src.Copy($null, $wb.Worksheets.Item(1))

# Freshly re-resolve the handle: any reference captured before a sheet
# insert/rename gets silently re-pointed to whatever now sits at that
# index, so we must grab it again right after the Copy() above.
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Re-resolve again after the rename, to be safe.
$q4 = $wb.Worksheets.Item(2)

# Grow from 10 rows (header + 9 data rows) to 14 rows (header + 13 data
# rows) by duplicating the format of an existing data row downward.
$q4.Range("A2:H2").Copy($q4.Range("A11:H11"))
$q4.Range("A2:H2").Copy($q4.Range("A12:H12"))
$q4.Range("A2:H2").Copy($q4.Range("A13:H13"))
$q4.Range("A2:H2").Copy($q4.Range("A14:H14"))

# Header row (unchanged text, already correct after the sheet Copy()).

function Set-FundRow {
    param($sheet, $row, $a, $code, $name, $scale, $pos, $ratio, $mktVal, $rank)

    $sheet.Range("A$row").Value = $a

    # Columns B..G must stay TEXT (fund codes have leading zeros, and the
    # template stores these as strings) - force text format, write the
    # value, then strip the style back to the sheet's default so no
    # stray numeric-format style id is introduced.
    $textRange = $sheet.Range("B$row`:F$row")
    $textRange.NumberFormat = "@"
    $sheet.Range("B$row").Value = $code
    $sheet.Range("C$row").Value = $name
    $sheet.Range("D$row").Value = $scale
    $sheet.Range("E$row").Value = $pos
    $sheet.Range("F$row").Value = $ratio
    $textRange.Style = "Normal"

    if ($mktVal -is [string]) {
        $gCell = $sheet.Range("G$row")
        $gCell.NumberFormat = "@"
        $gCell.Value = $mktVal
        $gCell.Style = "Normal"
    } else {
        $sheet.Range("G$row").Value = $mktVal
    }

    $sheet.Range("H$row").Value = $rank
}

Set-FundRow $q4 2  0  "290006" "泰信蓝筹精选混合"                 "15.07" "89.49" "2.55" "0.3843" 10
Set-FundRow $q4 3  1  "012850" "中融低碳经济3个月持有期混合A"       "5.54"  "88.38" "3.24" "0.1795" 8
Set-FundRow $q4 4  2  "001601" "鑫元鑫新收益灵活配置混合A"         "1.02"  "78.41" "3.83" "0.0391" 9
Set-FundRow $q4 5  3  "012851" "中融低碳经济3个月持有期混合C"       "1.19"  "88.38" "3.24" "0.0386" 8
Set-FundRow $q4 6  4  "014701" "中欧量化动能混合A"                "2.05"  "83.99" "1.11" "0.0228" 8
Set-FundRow $q4 7  5  "014702" "中欧量化动能混合C"                "1.93"  "83.99" "1.11" "0.0214" 8
Set-FundRow $q4 8  6  "161727" "招商增荣灵活配置混合（LOF）"       "0.52"  "55.51" "1.89" "0.0098" 9
Set-FundRow $q4 9  7  "005949" "鑫元行业轮动灵活配置混合A"         "0.18"  "67.42" "4.33" "0.0078" 6
Set-FundRow $q4 10 8  "009719" "招商增浩一年定期开放混合C"         "0.71"  "21.56" "0.98" "0.0070" 7
Set-FundRow $q4 11 9  "003670" "中融物联网主题灵活配置混合"         "0.13"  "92.35" "3.33" "0.0043" 6
Set-FundRow $q4 12 10 "009718" "招商增浩一年定期开放混合A"         "0.38"  "21.56" "0.98" "0.0037" 7
Set-FundRow $q4 13 11 "005950" "鑫元行业轮动灵活配置混合C"         "0.00"  "67.42" "4.33" 0        6
Set-FundRow $q4 14 12 "001602" "鑫元鑫新收益灵活配置混合C"         "0.00"  "78.41" "3.83" 0        9
